# Adds the "2022-Q4" quarter: a new summary row on "总计" plus a brand new
# "2022-Q4" worksheet (fund-holdings breakdown), inserted right after "总计"
# and before the existing "2022-Q3" tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for 2022-Q4 at the top of the
#    data block (row 2), pushing the existing quarters down by one row, and
#    renumber the helper index column (A).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A3").Copy() | Out-Null
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 24
$summary.Range("D2").Value = 3.19

# Renumber the index column for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# ---------------------------------------------------------------------------
# 2) New "2022-Q4" worksheet: fund-by-fund holdings breakdown, positioned
#    right after "总计" (i.e. before "2022-Q3"). Clone the "2022-Q3" sheet so
#    headers / column layout / styling come along for free, then overwrite
#    the data with the 2022-Q4 figures.
# ---------------------------------------------------------------------------
$fundData = @(
  @(0, "270021", "广发聚瑞混合A", "18.60", "90.06", "4.69", "0.8723", 7),
  @(1, "011136", "广发盛兴混合A", "17.77", "87.23", "3.91", "0.6948", 8),
  @(2, "506007", "广发科创板两年定开混合", "5.31", "88.81", "4.87", "0.2586", 6),
  @(3, "012342", "广发瑞泽精选混合A", "5.19", "89.55", "4.53", "0.2351", 6),
  @(4, "010161", "广发瑞安精选股票A", "6.49", "86.53", "3.58", "0.2323", 8),
  @(5, "006081", "海富通电子信息传媒产业股票A", "5.43", "92.10", "3.41", "0.1852", 8),
  @(6, "009623", "长城创新驱动混合A", "4.20", "91.95", "3.43", "0.1441", 7),
  @(7, "013000", "广发盛泽一年持有期混合A", "2.47", "85.08", "4.74", "0.1171", 6),
  @(8, "006080", "海富通电子信息传媒产业股票C", "3.10", "92.10", "3.41", "0.1057", 8),
  @(9, "011137", "广发盛兴混合C", "1.85", "87.23", "3.91", "0.0723", 8),
  @(10, "013346", "富荣信息技术混合C", "1.13", "91.06", "5.30", "0.0599", 5),
  @(11, "002133", "广发鑫益灵活配置混合", "1.33", "87.79", "4.04", "0.0537", 7),
  @(12, "004315", "前海开源沪港深新硬件主题灵活配置混合C", "0.91", "91.15", "3.35", "0.0305", 9),
  @(13, "013345", "富荣信息技术混合A", "0.45", "91.06", "5.30", "0.0238", 5),
  @(14, "004314", "前海开源沪港深新硬件主题灵活配置混合A", "0.68", "91.15", "3.35", "0.0228", 9),
  @(15, "010026", "广发聚瑞混合C", "0.44", "90.06", "4.69", "0.0206", 7),
  @(16, "010162", "广发瑞安精选股票C", "0.55", "86.53", "3.58", "0.0197", 8),
  @(17, "012343", "广发瑞泽精选混合C", "0.37", "89.55", "4.53", "0.0168", 6),
  @(18, "013001", "广发盛泽一年持有期混合C", "0.28", "85.08", "4.74", "0.0133", 6),
  @(19, "014598", "永赢合享混合A", "0.31", "31.53", "1.28", "0.0040", 7),
  @(20, "014433", "国泰智享科技1个月滚动持有混合A", "0.10", "61.54", "3.05", "0.0030", 10),
  @(21, "014599", "永赢合享混合C", "0.06", "31.53", "1.28", "0.0008", 7),
  @(22, "017458", "长城创新驱动混合C", "0.00", "91.95", "3.43", "__NUM0__", 7),
  @(23, "014434", "国泰智享科技1个月滚动持有混合C", "0.00", "61.54", "3.05", "__NUM0__", 10)
)

$q3 = $wb.Worksheets.Item("2022-Q3")
$insertBefore = $wb.Worksheets.Item(2)
$q3.Copy($insertBefore)
$new = $wb.Worksheets.Item(2)
$new.Name = "2022-Q4"

# The template ("2022-Q3") only has 12 rows (1 header + 11 data); extend the
# formatted block down to row 25 (1 header + 24 data rows) by cloning the
# format of the last template data row.
$new.Range("A12:H12").Copy() | Out-Null
$new.Range("A13:H25").PasteSpecial(-4122)

# Text-format columns B..G so numeric-looking strings (fund codes, scales,
# ratios, ...) are stored as text rather than being coerced to numbers -
# except G24:G25, which hold a genuine numeric 0.
$new.Range("B2:G23").NumberFormat = "@"
$new.Range("B24:F25").NumberFormat = "@"

$r = 2
foreach ($row in $fundData) {
    $new.Cells.Item($r, 1).Value = $row[0]
    $new.Cells.Item($r, 2).Value = $row[1]
    $new.Cells.Item($r, 3).Value = $row[2]
    $new.Cells.Item($r, 4).Value = $row[3]
    $new.Cells.Item($r, 5).Value = $row[4]
    $new.Cells.Item($r, 6).Value = $row[5]
    if ($row[6] -eq "__NUM0__") {
        $new.Cells.Item($r, 7).Value = 0
    } else {
        $new.Cells.Item($r, 7).Value = $row[6]
    }
    $new.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

Write-Output "done"
